# Apply "added subscription selection page" update to the Test Results sheet.
#
# Logical change (per the target diff):
#   - Row 2 (A2/C2): verifyCustomerNavigationToRegistrationPage / 29-10-2024
#                    -> verifyCustomerSuccessfulLogin / 03-11-2024
#   - Row 3 (A3/C3): verifyNewCustomerRegistrationSubmissionFlow / 29-10-2024
#                    -> verifyCustomerPreferredPackageSelection / 03-11-2024
#   - Row 4 (A4/C4): verifyCustomerRegistrationAndLoginNavigation / 29-10-2024
#                    -> verifyCustomerNavigationAfterSaving / 03-11-2024
#   - Rows 5 and 6 (verifyCustomerEmailActivation, verifyCustomerSuccessfulLogin)
#     are removed entirely, shrinking the table from 6 rows to 4 (header + 3).
#   - Column A is narrower afterwards because the longest remaining method name
#     is shorter than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Force the cell to stay a literal text value (type "s" in the OOXML,
    # i.e. a shared string) instead of letting Excel auto-convert
    # date-looking strings like "03-11-2024" into a numeric date serial.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Drop the temporary "Text" number format again so the cell ends up
    # with the same (default) style it had before, instead of keeping a
    # custom number-format style applied to it.
    $cell.ClearFormats()
}

# Row 2
$ws.Range("A2").Value = "verifyCustomerSuccessfulLogin"
Set-TextValue "C2" "03-11-2024"

# Row 3
$ws.Range("A3").Value = "verifyCustomerPreferredPackageSelection"
Set-TextValue "C3" "03-11-2024"

# Row 4
$ws.Range("A4").Value = "verifyCustomerNavigationAfterSaving"
Set-TextValue "C4" "03-11-2024"

# The two trailing rows (old verifyCustomerEmailActivation /
# verifyCustomerSuccessfulLogin entries) no longer exist in the refreshed
# report, so remove them completely (shifts dimension from A1:C6 to A1:C4).
$ws.Rows("5:6").Delete()

# Re-fit column A now that the longest method name in the sheet is shorter.
$ws.Columns("A").ColumnWidth = 37.91666666666667
